# Commit: "changed the look of transfer and added a button for adding account"
#
# This applies the underlying data/structural changes:
#  - Rename the two month sheets to proper capitalization (mars -> March, april -> April)
#  - On the March sheet, insert 12 new expense rows (rows 12-23) before the
#    "Monthly total" row, which now moves down to row 24 with an updated total.

$wb = $excel.ActiveWorkbook

$wsMarch = $wb.Worksheets.Item(1)
$wsApril = $wb.Worksheets.Item(2)

$wsMarch.Name = "March"
$wsApril.Name = "April"

# Insert 12 blank rows above the existing "Monthly total" row (currently row 12),
# pushing it down to row 24.
$wsMarch.Range("A12:A23").EntireRow.Insert()

# The Price/Date columns in this sheet store numeric- and date-looking values as
# plain text (shared strings), so force a text number format on the new block
# before writing values to keep Excel from auto-converting them to numbers/dates.
$wsMarch.Range("A12:E23").NumberFormat = "@"

$wsMarch.Range("A12").Value = "Food"
$wsMarch.Range("B12").Value = "pizza"
$wsMarch.Range("C12").Value = "2023-03-27"
$wsMarch.Range("D12").Value = "100.0"
$wsMarch.Range("E12").Value = "Checkings"

$wsMarch.Range("A13").Value = "Food"
$wsMarch.Range("B13").Value = "asfd"
$wsMarch.Range("C13").Value = "2023-03-27"
$wsMarch.Range("D13").Value = "100.0"
$wsMarch.Range("E13").Value = "Savings"

$wsMarch.Range("A14").Value = "Rent"
$wsMarch.Range("B14").Value = "dsf"
$wsMarch.Range("C14").Value = "2023-03-27"
$wsMarch.Range("D14").Value = "1000.0"
$wsMarch.Range("E14").Value = "Savings"

$wsMarch.Range("A15").Value = "Food"
$wsMarch.Range("B15").Value = "asdf"
$wsMarch.Range("C15").Value = "2023-03-27"
$wsMarch.Range("D15").Value = "1000.0"
$wsMarch.Range("E15").Value = "Savings"

$wsMarch.Range("A16").Value = "Food"
$wsMarch.Range("B16").Value = "elias"
$wsMarch.Range("C16").Value = "2023-03-27"
$wsMarch.Range("D16").Value = "1000.0"
$wsMarch.Range("E16").Value = "Savings"

$wsMarch.Range("A17").Value = "Food"
$wsMarch.Range("B17").Value = "banan"
$wsMarch.Range("C17").Value = "2023-03-27"
$wsMarch.Range("D17").Value = "1000.0"
$wsMarch.Range("E17").Value = "Savings"

$wsMarch.Range("A18").Value = "Transportation"
$wsMarch.Range("B18").Value = "buss"
$wsMarch.Range("C18").Value = "2023-03-27"
$wsMarch.Range("D18").Value = "1000.0"
$wsMarch.Range("E18").Value = "Savings"

$wsMarch.Range("A19").Value = "Transportation"
$wsMarch.Range("B19").Value = "buss"
$wsMarch.Range("C19").Value = "2023-03-27"
$wsMarch.Range("D19").Value = "1000.0"
$wsMarch.Range("E19").Value = "Checkings"

$wsMarch.Range("A20").Value = "Food"
$wsMarch.Range("B20").Value = "elias"
$wsMarch.Range("C20").Value = "2023-03-27"
$wsMarch.Range("D20").Value = "1000.0"
$wsMarch.Range("E20").Value = "Savings"

$wsMarch.Range("A21").Value = "Food"
$wsMarch.Range("B21").Value = "asf"
$wsMarch.Range("C21").Value = "2023-03-27"
$wsMarch.Range("D21").Value = "1000.0"
$wsMarch.Range("E21").Value = "Savings"

$wsMarch.Range("A22").Value = "Food"
$wsMarch.Range("B22").Value = "asf"
$wsMarch.Range("C22").Value = "2023-03-27"
$wsMarch.Range("D22").Value = "100.0"
$wsMarch.Range("E22").Value = "Card"

$wsMarch.Range("A23").Value = "Food"
$wsMarch.Range("B23").Value = "asdf"
$wsMarch.Range("C23").Value = "2023-03-27"
$wsMarch.Range("D23").Value = "234.0"
$wsMarch.Range("E23").Value = "Savings"

# Row 24 is the pre-existing "Monthly total" row, now shifted down from row 12.
# Its label stays textual, and its total is recalculated to include the new rows.
$wsMarch.Range("A24").Value = "Monthly total: "
$wsMarch.Range("B24").Value = 30894.0
